# Multi slots per student handling
# Rearranges / updates the per-slot Org<->Student assignment tables
# on the StudentsMapping and OrganizationMapping sheets.

$wb = $excel.ActiveWorkbook

$wsStudents = $wb.Worksheets.Item("StudentsMapping")

# Row 2
$wsStudents.Range("C2").Value = 'Keck VIO - COBI'
$wsStudents.Range("D2").Value = 'Cedars-Sinai - Neurosciences'
$wsStudents.Range("E2").Value = 'Providence Health Network'
$wsStudents.Range("F2").Value = 'Optum CF - Patient XP'

# Row 3
$wsStudents.Range("C3").Value = 'St.Johns-PhysPartners '
$wsStudents.Range("D3").Value = 'Optum CF - Patient XP'
$wsStudents.Range("E3").Value = 'Verdugo Hills Hospital'
$wsStudents.Range("F3").Value = 'Keck VIO - COBI'

# Row 4
$wsStudents.Range("C4").Value = 'CHLA - Anesthesia&CCM'
$wsStudents.Range("D4").Value = 'City of Hope - CMO'
$wsStudents.Range("E4").Value = 'Rancho Los Amigos NRC'
$wsStudents.Range("F4").Value = 'CHLA - Anesthesia&CCM'

# Row 5
$wsStudents.Range("C5").Value = 'City of Hope - CMO'
$wsStudents.Range("D5").Value = 'Emanate Health'
$wsStudents.Range("E5").Value = 'Cedars-Sinai - Neurosciences'
$wsStudents.Range("F5").Value = 'Rancho Los Amigos NRC'

# Row 6
$wsStudents.Range("C6").Value = 'Emanate Health'
$wsStudents.Range("D6").Value = 'Providence Health Network'
$wsStudents.Range("E6").Value = 'Keck VIO - COBI'
$wsStudents.Range("F6").Value = 'Kaiser PC - Consulting'

# Row 7
$wsStudents.Range("C7").Value = 'Kaiser PC - Consulting'
$wsStudents.Range("D7").Value = 'St.Johns-PhysPartners '
$wsStudents.Range("E7").Value = 'SCAN Health Plan'
$wsStudents.Range("F7").Value = 'Cedars-Sinai - Neurosciences'

# Row 8
$wsStudents.Range("C8").Value = 'Optum CF - Patient XP'
$wsStudents.Range("D8").Value = 'CHLA - Anesthesia&CCM'
$wsStudents.Range("E8").Value = 'City of Hope - CMO'
$wsStudents.Range("F8").Value = 'St.Johns-PhysPartners '

# Row 9
$wsStudents.Range("C9").Value = 'Cedars-Sinai - Neurosciences'
$wsStudents.Range("D9").Value = 'Keck VIO - COBI'
$wsStudents.Range("E9").Value = 'Optum CF - Patient XP'
$wsStudents.Range("F9").Value = 'Providence Health Network'

# Row 10
$wsStudents.Range("C10").Value = 'Keck IRM'
$wsStudents.Range("D10").Value = 'West Hills Hospital'
$wsStudents.Range("E10").Value = 'Kaiser PC - Consulting'
$wsStudents.Range("F10").Value = 'Verdugo Hills Hospital'

# Row 11
$wsStudents.Range("C11").Value = 'Torrance Memorial'
$wsStudents.Range("D11").Value = 'Optum CF - Digi Transformation'
$wsStudents.Range("E11").Value = 'West Hills Hospital'
$wsStudents.Range("F11").Value = 'Keck IRM'

# Row 12
$wsStudents.Range("C12").Value = 'Verdugo Hills Hospital'
$wsStudents.Range("D12").Value = 'Kaiser PC - Consulting'
$wsStudents.Range("E12").Value = 'Optum CF - Digi Transformation'
$wsStudents.Range("F12").Value = 'City of Hope - CMO'

# Row 13
$wsStudents.Range("C13").Value = 'West Hills Hospital'
$wsStudents.Range("D13").Value = 'Verdugo Hills Hospital'
$wsStudents.Range("E13").Value = 'CHLA - Anesthesia&CCM'
$wsStudents.Range("F13").Value = 'Optum CF - Digi Transformation'

$wsOrg = $wb.Worksheets.Item("OrganizationMapping")

# Row 2
$wsOrg.Range("C2").Value = 'Raashi Subramanya'
$wsOrg.Range("D2").Value = 'Oceana Hanner'
$wsOrg.Range("E2").Value = 'Fahima Gohil'
$wsOrg.Range("F2").Value = 'Emma Crusinberry'

# Row 3
$wsOrg.Range("C3").Value = 'Daniela Ahumada'
$wsOrg.Range("D3").Value = 'Stanley Ibe'
$wsOrg.Range("E3").Value = ' Bryce Dechert'
$wsOrg.Range("F3").Value = 'Daniela Ahumada'

# Row 4
$wsOrg.Range("C4").Value = 'Fahima Gohil'
$wsOrg.Range("D4").Value = 'Daniela Ahumada'
$wsOrg.Range("E4").Value = 'Stanley Ibe'
$wsOrg.Range("F4").Value = 'Russelle Chang'

# Row 5
$wsOrg.Range("C5").Value = 'Julia Orozco'
$wsOrg.Range("D5").Value = 'Fahima Gohil'
$wsOrg.Range("E5").Value = ""
$wsOrg.Range("F5").Value = ""

# Row 6
$wsOrg.Range("C6").Value = 'Emma Crusinberry'
$wsOrg.Range("D6").Value = 'Russelle Chang'
$wsOrg.Range("E6").Value = 'Nikhil Bajpai'
$wsOrg.Range("F6").Value = 'Julia Orozco'

# Row 7
$wsOrg.Range("C7").Value = 'Nikhil Bajpai'
$wsOrg.Range("D7").Value = ""
$wsOrg.Range("E7").Value = ""
$wsOrg.Range("F7").Value = 'Eryn Burnette'

# Row 8
$wsOrg.Range("C8").Value = 'Oceana Hanner'
$wsOrg.Range("D8").Value = 'Raashi Subramanya'
$wsOrg.Range("E8").Value = 'Julia Orozco'
$wsOrg.Range("F8").Value = 'Esther Choi'

# Row 9
$wsOrg.Range("C9").Value = ""
$wsOrg.Range("D9").Value = 'Eryn Burnette'
$wsOrg.Range("E9").Value = 'Russelle Chang'
$wsOrg.Range("F9").Value = ' Bryce Dechert'

# Row 10
$wsOrg.Range("C10").Value = 'Stanley Ibe'
$wsOrg.Range("D10").Value = 'Esther Choi'
$wsOrg.Range("E10").Value = 'Raashi Subramanya'
$wsOrg.Range("F10").Value = 'Oceana Hanner'

# Row 11
$wsOrg.Range("C11").Value = ""
$wsOrg.Range("D11").Value = 'Julia Orozco'
$wsOrg.Range("E11").Value = 'Oceana Hanner'
$wsOrg.Range("F11").Value = 'Raashi Subramanya'

# Row 12
$wsOrg.Range("C12").Value = ""
$wsOrg.Range("D12").Value = ""
$wsOrg.Range("E12").Value = 'Daniela Ahumada'
$wsOrg.Range("F12").Value = 'Fahima Gohil'

# Row 13
$wsOrg.Range("C13").Value = ""
$wsOrg.Range("D13").Value = ""
$wsOrg.Range("E13").Value = 'Emma Crusinberry'
$wsOrg.Range("F13").Value = ""

# Row 14
$wsOrg.Range("C14").Value = 'Esther Choi'
$wsOrg.Range("D14").Value = 'Emma Crusinberry'
$wsOrg.Range("E14").Value = ""
$wsOrg.Range("F14").Value = 'Stanley Ibe'

# Row 15
$wsOrg.Range("C15").Value = 'Eryn Burnette'
$wsOrg.Range("D15").Value = ""
$wsOrg.Range("E15").Value = ""
$wsOrg.Range("F15").Value = ""

# Row 16
$wsOrg.Range("C16").Value = 'Russelle Chang'
$wsOrg.Range("D16").Value = ' Bryce Dechert'
$wsOrg.Range("E16").Value = 'Esther Choi'
$wsOrg.Range("F16").Value = 'Nikhil Bajpai'

# Row 17
$wsOrg.Range("C17").Value = ' Bryce Dechert'
$wsOrg.Range("D17").Value = 'Nikhil Bajpai'
$wsOrg.Range("E17").Value = 'Eryn Burnette'
$wsOrg.Range("F17").Value = ""

# New (empty) cell in column G - extends the sheet's used range to A1:G17
$wsOrg.Range("G7").Font.Bold = $wsOrg.Range("G7").Font.Bold
